$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28 - this shifts existing rows 28-32 down to 29-33
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new weekly record
$ws.Range("A28").Value = 6
$ws.Range("B28").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 44776
$ws.Range("D28").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 100112035
$ws.Range("G28").Value = "Bruselas (repollito)"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 580
$ws.Range("K28").Value = 17000
$ws.Range("L28").Value = 19000
$ws.Range("M28").Value = 17897
$ws.Range("N28").Value = "$/malla 15 kilos"
$ws.Range("O28").Value = "Provincia de Quillota"
$ws.Range("P28").Value = 1193
$ws.Range("Q28").Value = 15
$ws.Range("R28").Value = "Hortaliza"
